$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # METADATA
$ws2 = $wb.Worksheets.Item(2)   # CONFIGURATION

# --- METADATA sheet: insert a new column F for "CATALOG_IDENTIFIER" / "catalogIdentifier" ---
# (shifts old F..I -> G..J). This also drops the now-pointless empty placeholder
# cells in the newly vacated column F for rows 2-7 (bugfix for completely empty rows/cols).
$ws1.Columns("F").Insert() | Out-Null
$ws1.Range("F1").Value = "CATALOG_IDENTIFIER"
$ws1.Columns("F").ColumnWidth = 20

# --- CONFIGURATION sheet: add row 11 describing the new bean property mapping ---
$ws2.Range("A11").Value = "CATALOG_IDENTIFIER"
$ws2.Range("C11").Value = "catalogIdentifier"

# --- Update selections / active sheet/tab ---
$ws1.Range("F1").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B13").Select() | Out-Null
